$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBS")

# Append a new row (row 9) of lookup metadata, matching the pattern of the
# preceding rows (Key ID / 讀取Key條件 / 其他ORDER條件 triplet).
$ws.Range("A9").Value = "findCustNoAndActualRepayDateFirst"
$ws.Range("B9").Value = "CustNo = ,AND ActualRepayDate >= ,AND ActualRepayDate <="
$ws.Range("C9").Value = "ActualRepayDate DESC,LogNo DESC"

# Match the saved selection state (active cell moved to A9 on that sheet).
$ws.Activate()
$ws.Range("A9").Select()
